$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Dropdown Example"

# Header cell
$ws.Range("A1").Value = "Select Option"

# Touch A2 (with the default "Normal" style) so the sheet's used range
# extends to A2, matching the helper-row added alongside the new validation.
$ws.Range("A2").Style = "Normal"

# Replace the old fruit dropdown (B2:B10) with a construction-progress
# dropdown on A2.
$ws.Range("B2:B10").Validation.Delete()
$ws.Range("A2").Validation.Add(3, 1, 1, '"0 - 10%: Foundation completed: Groundwork finished; no vertical structure yet.,11 - 25%: Structure and rough-in started: Structural framing in progress; initial MEP rough-in."')
$ws.Range("A2").Validation.IgnoreBlank = $true
$ws.Range("A2").Validation.InCellDropdown = $true
$ws.Range("A2").Validation.ShowInput = $false
$ws.Range("A2").Validation.ShowError = $false
